$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the bold/bordered
# header style from H1 so the new headers match the existing ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I and J columns (rows 2-70) with their data values.
$iValues = @(8,9,8,9,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,7,9,7,9,9,9,7,8,8,9,8,9,9,9,8,9,9,8,8,9,8,9,7,7,9,9,7,8,8,9,9,9,8,7,10,9,8,8,8,7,8,9,6,5,6,7,9,7,7)
$jValues = @(9,9,9,9,9,9,9,7,9,9,9,9,9,9,10,9,9,9,9,9,7,9,7,9,9,9,7,8,8,9,8,9,9,9,9,9,9,9,8,9,8,9,7,7,9,9,7,8,9,9,9,9,8,7,10,9,8,8,8,7,8,9,6,5,6,7,9,7,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
